$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns, copying the formatting
# (bold font + border + center/top alignment) from the existing
# header cell AB1, then overwrite the displayed text.
$ws.Range("AB1").Copy($ws.Range("AC1"))
$ws.Range("AB1").Copy($ws.Range("AD1"))
$ws.Range("AB1").Copy($ws.Range("AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 77
    $ws.Cells.Item($r, 30).Value = 83
    $ws.Cells.Item($r, 31).Value = 0
}
